# Macroferia Regional de Talca - Zapallo italiano
# A new weekly price record is inserted as row 181; all the existing
# records that were previously rows 181-195 shift down to rows 182-196.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 181 (pushes 181:195 -> 182:196, extends the used
# range to row 196, and carries the date-format style from column D down
# into the freshly inserted row, same as Excel's native Insert behaviour).
$ws.Rows("181:181").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A181").Value = 5
$ws.Range("B181").Value = "Macroferia Regional de Talca"
$ws.Range("C181").Value = "Maule"
$ws.Range("D181").Value = 44461
$ws.Range("E181").Value = 7
$ws.Range("F181").Value = 100112032
$ws.Range("G181").Value = "Zapallo italiano"
$ws.Range("H181").Value = "Sin especificar"
$ws.Range("I181").Value = "Primera"
$ws.Range("J181").Value = 400
$ws.Range("K181").Value = 9000
$ws.Range("L181").Value = 9000
$ws.Range("M181").Value = 9000
$ws.Range("N181").Value = "`$/caja 50 unidades"
$ws.Range("O181").Value = "Región de Arica y Parinacota"
$ws.Range("P181").Value = 180
$ws.Range("Q181").Value = 50
$ws.Range("R181").Value = "Hortaliza"
